$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 276
$ws.Range("F5").Value = 2901
$ws.Range("F8").Value = 2266
$ws.Range("F9").Value = 1503
$ws.Range("F10").Value = 43
$ws.Range("F13").Value = 2597
$ws.Range("F16").Value = 6096
$ws.Range("F18").Value = 5488
$ws.Range("F19").Value = 3
$ws.Range("F20").Value = 2094
$ws.Range("F21").Value = 2975
$ws.Range("F22").Value = 3402
$ws.Range("F23").Value = 199
$ws.Range("F24").Value = 1671
$ws.Range("F25").Value = 37
$ws.Range("F28").Value = 153
$ws.Range("F30").Value = 346
$ws.Range("F32").Value = 2262
$ws.Range("F34").Value = 137
$ws.Range("F35").Value = 321
$ws.Range("F36").Value = 853
$ws.Range("F38").Value = 406
$ws.Range("F39").Value = 471
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 29
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 276
$ws.Range("F9").Value = 2901
$ws.Range("F11").Value = 2266
$ws.Range("F12").Value = 1503
$ws.Range("F13").Value = 43
$ws.Range("F17").Value = 2597
$ws.Range("F22").Value = 29
$ws.Range("F23").Value = 6096
$ws.Range("F25").Value = 5488
$ws.Range("F26").Value = 2094
$ws.Range("F27").Value = 2975
$ws.Range("F28").Value = 3402
$ws.Range("F30").Value = 199
$ws.Range("F33").Value = 1671
$ws.Range("F38").Value = 153
$ws.Range("F40").Value = 346
$ws.Range("F42").Value = 2262
$ws.Range("F44").Value = 137
$ws.Range("F45").Value = 321
$ws.Range("F46").Value = 853
$ws.Range("F48").Value = 406
$ws.Range("F49").Value = 471
